# Reorders the header columns (A1:F1) and updates the corresponding
# block-order indicator values (rows 2-7) to match the new column order.
#
# New header order:  kitchens_1 | bedrooms_1 | bedrooms_2 | kitchens_2 | living_rooms_1 | living_rooms_2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers ---
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "bedrooms_2"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# --- Data rows 2-7 (0/1 indicator matrix), only cells that actually changed ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 1

$ws.Range("A3").Value = 1
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = 0
$ws.Range("F4").Value = 1

$ws.Range("D6").Value = 1
$ws.Range("F6").Value = 0
